# Fruta / hortaliza, semanal
# Insert 3 new weekly price rows above the existing data block (new rows 536-538),
# pushing the previous rows 536-553 down to 556.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three whole rows at row 536 (existing rows 536:553 shift down to 539:556)
$ws.Range("536:538").Insert()

$newRows = @(
    @(5, "Macroferia Regional de Talca", "Maule", 44939, 7, "Fruta", 100109, "Uva", 100109001, "Uva", "Flame Seedless", "Primera", 230, 10000, 10000, 10000, "`$/bandeja 10 kilos", "Provincia de Limarí", 1000, 10),
    @(5, "Macroferia Regional de Talca", "Maule", 44939, 7, "Fruta", 100109, "Uva", 100109001, "Uva", "Red Globe", "Primera", 180, 12000, 12000, 12000, "`$/bandeja 10 kilos", "Provincia de Limarí", 1200, 10),
    @(5, "Macroferia Regional de Talca", "Maule", 44939, 7, "Fruta", 100109, "Uva", 100109001, "Uva", "Superior Seedless", "Primera", 200, 10000, 10000, 10000, "`$/bandeja 10 kilos", "Provincia de Limarí", 1000, 10)
)

$r = 536
foreach ($row in $newRows) {
    for ($i = 0; $i -lt $row.Length; $i++) {
        $ws.Cells.Item($r, $i + 1).Value2 = $row[$i]
    }
    $r = $r + 1
}
